# helping-game/tableconfiguration.xlsx
# - add a thin right border to the "group boundary" columns (D and H) for
#   rows 1-4, keeping the existing center alignment
# - replace the three generic "white" placeholder cells in row 4 with
#   proper per-group labels (A4, B4, C4)
# - move the active selection from F6 to D1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and H mark the right edge of each 4-column group (A-D, E-H, I-L).
# Give those cells a thin right border in addition to their existing
# center alignment.
$boundaryCells = @("D1", "H1", "D2", "H2", "D3", "H3", "D4", "H4")
foreach ($cellRef in $boundaryCells) {
    $cell = $ws.Range($cellRef)
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight = thin continuous line
}

# Row 4 used a shared "white" placeholder string for the per-group header
# cells; replace each with a distinct label matching its column group.
$ws.Range("D4").Value = "A4"
$ws.Range("H4").Value = "B4"
$ws.Range("L4").Value = "C4"

# Update the sheet's active selection.
$ws.Range("D1").Select()
